# Actualizacion de la estimacion - 4ta iteracion de la fase de construccion
$wb = $excel.ActiveWorkbook

$pcu    = $wb.Worksheets.Item("PCU")
$casos  = $wb.Worksheets.Item("Casos De Uso")
$actores = $wb.Worksheets.Item("Actores")
$hh     = $wb.Worksheets.Item("Estimacion h-h")

# Tabla "Caracteristicas especiales" (factor tecnico) - valores revisados
$pcu.Range("F24").Value = 0   # antes 3 - Concurrencia
$pcu.Range("F27").Value = 0   # antes 1 - Facilidad de instalacion

# Tabla "Factores de entorno" - valores revisados
$pcu.Range("G33").Value = 5   # antes 4 - Familiaridad con un proceso definido
$pcu.Range("G39").Value = 5   # antes 3 - Requerimientos estables
$pcu.Range("G40").Value = 1   # antes 2 - Miembros a tiempo parcial

$excel.Calculate()

# Posiciones de seleccion finales tal cual quedaron en cada hoja
$pcu.Range("G24").Select()
$casos.Range("C35").Select()
$actores.Range("F15").Select()
$hh.Range("D23").Select()
